$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change the header of column B from "Count" to "House Sold"
$ws.Range("B1").Value = "House Sold"

# Update the active selection to D4, matching the diff
$ws.Range("D4").Select()
